$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shuffled match rows (home/away swapped in source feed) ---
# Row 12
$ws.Range("F12").Value = "Mamelodi Sundowns"
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = "Kaizer Chiefs"
$ws.Range("J12").Value = 1.77
$ws.Range("K12").Value = "06/08/2023 17:12"
$ws.Range("L12").Value = 1.45
$ws.Range("M12").Value = "09/08/2023 14:59"
$ws.Range("N12").Value = 3.45
$ws.Range("O12").Value = "06/08/2023 17:12"
$ws.Range("P12").Value = 4.3
$ws.Range("Q12").Value = "09/08/2023 14:59"
$ws.Range("R12").Value = 4.96
$ws.Range("S12").Value = "06/08/2023 17:12"
$ws.Range("T12").Value = 7.82
$ws.Range("U12").Value = "09/08/2023 14:59"
$ws.Range("V12").Value = "https://www.betexplorer.com/football/south-africa/premier-league/mamelodi-sundowns-kaizer-chiefs/GCKwva1T/"

# Row 13
$ws.Range("F13").Value = "Swallows"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = "Cape Town City"
$ws.Range("J13").Value = 2.88
$ws.Range("K13").Value = "06/08/2023 17:42"
$ws.Range("L13").Value = 3.03
$ws.Range("M13").Value = "09/08/2023 14:50"
$ws.Range("N13").Value = 2.81
$ws.Range("O13").Value = "06/08/2023 17:42"
$ws.Range("P13").Value = 2.72
$ws.Range("Q13").Value = "09/08/2023 14:50"
$ws.Range("R13").Value = 2.86
$ws.Range("S13").Value = "06/08/2023 17:42"
$ws.Range("T13").Value = 2.9
$ws.Range("U13").Value = "09/08/2023 14:50"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/south-africa/premier-league/swallows-fc-cape-town-city/ro7WuLWG/"

# Row 15
$ws.Range("F15").Value = "Polokwane"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "Stellenbosch"
$ws.Range("J15").Value = 3.15
$ws.Range("K15").Value = "08/08/2023 13:42"
$ws.Range("L15").Value = 3.27
$ws.Range("M15").Value = "09/08/2023 19:21"
$ws.Range("N15").Value = 3.04
$ws.Range("O15").Value = "08/08/2023 13:42"
$ws.Range("P15").Value = 2.85
$ws.Range("Q15").Value = "09/08/2023 19:21"
$ws.Range("R15").Value = 2.46
$ws.Range("S15").Value = "08/08/2023 13:42"
$ws.Range("T15").Value = 2.58
$ws.Range("U15").Value = "09/08/2023 19:21"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/south-africa/premier-league/polokwane-city-stellenbosch-fc/l0khcM1i/"

# Row 16
$ws.Range("F16").Value = "Chippa Utd."
$ws.Range("H16").Value = "TS Galaxy"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2.61
$ws.Range("K16").Value = "06/08/2023 21:12"
$ws.Range("L16").Value = 2.68
$ws.Range("M16").Value = "09/08/2023 19:20"
$ws.Range("N16").Value = 2.95
$ws.Range("O16").Value = "06/08/2023 21:12"
$ws.Range("P16").Value = 2.76
$ws.Range("Q16").Value = "09/08/2023 19:25"
$ws.Range("R16").Value = 3.1
$ws.Range("S16").Value = "06/08/2023 21:12"
$ws.Range("T16").Value = 3.24
$ws.Range("U16").Value = "09/08/2023 19:20"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/south-africa/premier-league/chippa-utd-ts-galaxy/zc1AYxPj/"

# Row 17
$ws.Range("F17").Value = "Cape Town Spurs"
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = "Sekhukhune"
$ws.Range("I17").Value = 2
$ws.Range("K17").Value = "07/08/2023 22:12"
$ws.Range("L17").Value = 2.9
$ws.Range("M17").Value = "09/08/2023 19:26"
$ws.Range("N17").Value = 2.96
$ws.Range("O17").Value = "07/08/2023 22:12"
$ws.Range("P17").Value = 2.77
$ws.Range("Q17").Value = "09/08/2023 19:26"
$ws.Range("R17").Value = 2.51
$ws.Range("S17").Value = "07/08/2023 22:12"
$ws.Range("T17").Value = 2.96
$ws.Range("U17").Value = "09/08/2023 19:06"
$ws.Range("V17").Value = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-sekhukhune/SjlddtHc/"

# Row 37
$ws.Range("F37").Value = "AmaZulu"
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = "TS Galaxy"
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 2.58
$ws.Range("K37").Value = "28/08/2023 13:23"
$ws.Range("L37").Value = 3.23
$ws.Range("M37").Value = "30/08/2023 19:19"
$ws.Range("N37").Value = 2.77
$ws.Range("O37").Value = "28/08/2023 13:23"
$ws.Range("P37").Value = 2.75
$ws.Range("Q37").Value = "30/08/2023 19:19"
$ws.Range("R37").Value = 3.1
$ws.Range("S37").Value = "28/08/2023 13:23"
$ws.Range("T37").Value = 2.7
$ws.Range("U37").Value = "30/08/2023 19:19"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/south-africa/premier-league/amazulu-ts-galaxy/0xOg1Sbf/"

# Row 38
$ws.Range("F38").Value = "Richards Bay"
$ws.Range("H38").Value = "Sekhukhune"
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 2.96
$ws.Range("K38").Value = "28/08/2023 13:24"
$ws.Range("L38").Value = 3.28
$ws.Range("M38").Value = "30/08/2023 19:29"
$ws.Range("N38").Value = 2.79
$ws.Range("O38").Value = "28/08/2023 13:24"
$ws.Range("P38").Value = 2.6
$ws.Range("Q38").Value = "30/08/2023 19:29"
$ws.Range("R38").Value = 2.67
$ws.Range("S38").Value = "28/08/2023 13:24"
$ws.Range("T38").Value = 2.82
$ws.Range("U38").Value = "30/08/2023 19:29"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-sekhukhune/djxUp7ED/"

# Row 39
$ws.Range("F39").Value = "Chippa Utd."
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = "Royal AM"
$ws.Range("I39").Value = 3
$ws.Range("J39").Value = 2.24
$ws.Range("L39").Value = 2.16
$ws.Range("M39").Value = "30/08/2023 19:27"
$ws.Range("N39").Value = 3
$ws.Range("P39").Value = 3.01
$ws.Range("Q39").Value = "30/08/2023 19:27"
$ws.Range("R39").Value = 3.42
$ws.Range("T39").Value = 4
$ws.Range("U39").Value = "30/08/2023 19:27"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/south-africa/premier-league/chippa-utd-royal-am/SOPo3lUs/"

# --- Append new match rows 86:87 (copy formatting from last existing row) ---
$ws.Range("A85:V85").Copy()
$ws.Range("A86:V87").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 86
$ws.Range("A86").Value = 85
$ws.Range("B86").Value = "south-africa"
$ws.Range("C86").Value = "premier-league"
$ws.Range("D86").Value = "2023-2024"
$ws.Range("E86").Value = 45242.60416666666
$ws.Range("F86").Value = "Golden Arrows"
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = "AmaZulu"
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 2.31
$ws.Range("K86").Value = "11/11/2023 13:45"
$ws.Range("L86").Value = 2.67
$ws.Range("M86").Value = "12/11/2023 14:28"
$ws.Range("N86").Value = 2.92
$ws.Range("O86").Value = "11/11/2023 13:45"
$ws.Range("P86").Value = 2.85
$ws.Range("Q86").Value = "12/11/2023 14:26"
$ws.Range("R86").Value = 3.38
$ws.Range("S86").Value = "11/11/2023 13:45"
$ws.Range("T86").Value = 3.14
$ws.Range("U86").Value = "12/11/2023 14:28"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/south-africa/premier-league/golden-arrows-amazulu/zgmhTd51/"

# Row 87
$ws.Range("A87").Value = 86
$ws.Range("B87").Value = "south-africa"
$ws.Range("C87").Value = "premier-league"
$ws.Range("D87").Value = "2023-2024"
$ws.Range("E87").Value = 45242.60416666666
$ws.Range("F87").Value = "Richards Bay"
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = "Cape Town Spurs"
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 2.48
$ws.Range("K87").Value = "11/11/2023 13:45"
$ws.Range("L87").Value = 2.03
$ws.Range("M87").Value = "12/11/2023 14:21"
$ws.Range("N87").Value = 2.88
$ws.Range("O87").Value = "11/11/2023 13:45"
$ws.Range("P87").Value = 3.15
$ws.Range("Q87").Value = "12/11/2023 14:21"
$ws.Range("R87").Value = 3.13
$ws.Range("S87").Value = "11/11/2023 13:45"
$ws.Range("T87").Value = 4.22
$ws.Range("U87").Value = "12/11/2023 14:21"
$ws.Range("V87").Value = "https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-cape-town-spurs/EqndSGK7/"

[void]$ws.Range("A1").Select()
